# Updated data to reflect new requirement separation
# The old column D ("Terms Typically Offered") is split so that three new
# columns (Corequisites, Concurrent, Recommended) are inserted before it.
# Every existing "Terms Typically Offered" value moves from D to the new G
# column, and the new D/E/F columns are populated with "NA" for each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---------------------------------------------------
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# --- Data rows (2-26): move old "Terms Typically Offered" value (col D) to
#     the new col G, then fill D/E/F with "NA" ----------------------------
for ($row = 2; $row -le 26; $row++) {
    $oldD = $ws.Cells.Item($row, 4).Text
    $ws.Cells.Item($row, 7).Value = $oldD
    $ws.Cells.Item($row, 4).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
}
